$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update product class (itemsTypeVerification) from Electronics to Mobiles
$ws.Range("B6").Value = "Mobiles"

# Add new row for product name / class ("login" test data)
$ws.Range("B7").Value = "Men"
$ws.Range("A7").Value = "productName"

# Update login / view state: ruler visible again, new selection
$excel.ActiveWindow.DisplayRuler = $true
$ws.Range("B11").Select()
